$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(49, 8).Value = 1309.4
$ws.Cells.Item(49, 10).Value = 2082.3333
$ws.Cells.Item(49, 12).Value = 6246.999899999999
$ws.Cells.Item(49, 14).Value = -6518.999899999999

$ws.Cells.Item(80, 8).Value = 645
$ws.Cells.Item(80, 9).Value = 306
$ws.Cells.Item(80, 10).Value = 899.25
$ws.Cells.Item(80, 11).Value = 918
$ws.Cells.Item(80, 12).Value = 2697.75
$ws.Cells.Item(80, 13).Value = 80
$ws.Cells.Item(80, 14).Value = -4693.75

$ws.Cells.Item(83, 8).Value = 645
$ws.Cells.Item(83, 9).Value = 306
$ws.Cells.Item(83, 10).Value = 899.25
$ws.Cells.Item(83, 11).Value = 2754
$ws.Cells.Item(83, 12).Value = 8093.25
$ws.Cells.Item(83, 13).Value = 2238
$ws.Cells.Item(83, 14).Value = -18077.25

$ws.Cells.Item(86, 8).Value = 1895.5454
$ws.Cells.Item(86, 9).Value = 1640.6
$ws.Cells.Item(86, 11).Value = 1640.6
$ws.Cells.Item(86, 13).Value = -517.5999999999999

$ws.Cells.Item(89, 8).Value = 1895.5454
$ws.Cells.Item(89, 9).Value = 1640.6
$ws.Cells.Item(89, 11).Value = 8203
$ws.Cells.Item(89, 13).Value = -2587

$ws.Cells.Item(106, 8).Value = 12108.2
$ws.Cells.Item(106, 9).Value = 1818
$ws.Cells.Item(106, 11).Value = 1818
$ws.Cells.Item(106, 13).Value = -1187

$ws.Cells.Item(112, 8).Value = 3829.4736
$ws.Cells.Item(112, 9).Value = 2666.5
$ws.Cells.Item(112, 10).Value = 4139.6
$ws.Cells.Item(112, 11).Value = 7999.5
$ws.Cells.Item(112, 12).Value = 12418.8
$ws.Cells.Item(112, 13).Value = -6891.5
$ws.Cells.Item(112, 14).Value = -14634.8

$ws.Cells.Item(137, 8).Value = 4235.4165
$ws.Cells.Item(137, 9).Value = 1466
$ws.Cells.Item(137, 11).Value = 4398
$ws.Cells.Item(137, 13).Value = -1848

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 53073.17
$ws.Cells.Item(32, 9).Value = 30470.172
$ws.Cells.Item(32, 11).Value = 30470.172
$ws.Cells.Item(32, 13).Value = -30183.172

$ws.Cells.Item(46, 8).Value = 10010
$ws.Cells.Item(46, 10).Value = 10010
$ws.Cells.Item(46, 12).Value = 10010
$ws.Cells.Item(46, 14).Value = -10648

$ws.Cells.Item(74, 8).Value = 1088.0714
$ws.Cells.Item(74, 9).Value = 1002.7619
$ws.Cells.Item(74, 11).Value = 1002.7619
$ws.Cells.Item(74, 13).Value = -128.7619

$ws.Cells.Item(77, 8).Value = 1088.0714
$ws.Cells.Item(77, 9).Value = 1002.7619
$ws.Cells.Item(77, 11).Value = 5013.809499999999
$ws.Cells.Item(77, 13).Value = -645.8094999999994

$ws.Cells.Item(102, 8).Value = 2703.3076
$ws.Cells.Item(102, 9).Value = 2589.875
$ws.Cells.Item(102, 11).Value = 2589.875
$ws.Cells.Item(102, 13).Value = -967.875

$ws.Cells.Item(122, 8).Value = 1640.826
$ws.Cells.Item(122, 9).Value = 1463.762
$ws.Cells.Item(122, 11).Value = 4391.286
$ws.Cells.Item(122, 13).Value = -1941.286

$ws.Cells.Item(132, 8).Value = 1379.2941
$ws.Cells.Item(132, 9).Value = 1059.1464
$ws.Cells.Item(132, 11).Value = 3177.4392
$ws.Cells.Item(132, 13).Value = -647.4392000000003

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 10086.529
$ws.Cells.Item(20, 9).Value = 9361
$ws.Cells.Item(20, 11).Value = 9361
$ws.Cells.Item(20, 13).Value = -9114

$ws.Cells.Item(134, 8).Value = 1159.5238
$ws.Cells.Item(134, 9).Value = 1152.5
$ws.Cells.Item(134, 11).Value = 3457.5
$ws.Cells.Item(134, 13).Value = -922.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(7, 8).Value = 28571750
$ws.Cells.Item(7, 9).Value = 50000252
$ws.Cells.Item(7, 10).Value = 417.86667
$ws.Cells.Item(7, 11).Value = 50000252
$ws.Cells.Item(7, 12).Value = 417.86667
$ws.Cells.Item(7, 13).Value = -50000139
$ws.Cells.Item(7, 14).Value = -643.86667

$ws.Cells.Item(58, 8).Value = 1096.9166
$ws.Cells.Item(58, 9).Value = 1008.5
$ws.Cells.Item(58, 11).Value = 1008.5
$ws.Cells.Item(58, 13).Value = -805.5

$ws.Cells.Item(134, 8).Value = 2384.2173
$ws.Cells.Item(134, 9).Value = 2547.25
$ws.Cells.Item(134, 11).Value = 7641.75
$ws.Cells.Item(134, 13).Value = -5106.75

$ws.Cells.Item(136, 8).Value = 1096.9166
$ws.Cells.Item(136, 9).Value = 1008.5
$ws.Cells.Item(136, 11).Value = 3025.5
$ws.Cells.Item(136, 13).Value = -475.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(50, 8).Value = 200146.4
$ws.Cells.Item(50, 10).Value = 0
$ws.Cells.Item(50, 12).Value = 0
$ws.Cells.Item(50, 14).ClearContents()

$ws.Cells.Item(53, 8).Value = 200146.4
$ws.Cells.Item(53, 10).Value = 0
$ws.Cells.Item(53, 12).Value = 0
$ws.Cells.Item(53, 14).ClearContents()

$ws.Cells.Item(63, 8).Value = 0
$ws.Cells.Item(63, 9).Value = 0
$ws.Cells.Item(63, 11).Value = 0
$ws.Cells.Item(63, 13).ClearContents()

$ws.Cells.Item(66, 8).Value = 0
$ws.Cells.Item(66, 9).Value = 0
$ws.Cells.Item(66, 11).Value = 0
$ws.Cells.Item(66, 13).ClearContents()

$ws.Cells.Item(69, 8).Value = 4598
$ws.Cells.Item(69, 9).Value = 6331.3335
$ws.Cells.Item(69, 10).Value = 1998
$ws.Cells.Item(69, 11).Value = 18994.0005
$ws.Cells.Item(69, 12).Value = 5994
$ws.Cells.Item(69, 13).Value = -18183.0005
$ws.Cells.Item(69, 14).Value = -7616

$ws.Cells.Item(72, 8).Value = 4598
$ws.Cells.Item(72, 9).Value = 6331.3335
$ws.Cells.Item(72, 10).Value = 1998
$ws.Cells.Item(72, 11).Value = 56982.0015
$ws.Cells.Item(72, 12).Value = 17982
$ws.Cells.Item(72, 13).Value = -52926.0015
$ws.Cells.Item(72, 14).Value = -26094

$ws.Cells.Item(74, 8).Value = 0
$ws.Cells.Item(74, 10).Value = 0
$ws.Cells.Item(74, 12).Value = 0
$ws.Cells.Item(74, 14).ClearContents()

$ws.Cells.Item(77, 8).Value = 0
$ws.Cells.Item(77, 10).Value = 0
$ws.Cells.Item(77, 12).Value = 0
$ws.Cells.Item(77, 14).ClearContents()

$ws.Cells.Item(80, 8).Value = 0
$ws.Cells.Item(80, 10).Value = 0
$ws.Cells.Item(80, 12).Value = 0
$ws.Cells.Item(80, 14).ClearContents()

$ws.Cells.Item(81, 8).Value = 100017300
$ws.Cells.Item(81, 10).Value = 100017300
$ws.Cells.Item(81, 12).Value = 300051900
$ws.Cells.Item(81, 14).Value = -300054146

$ws.Cells.Item(83, 8).Value = 0
$ws.Cells.Item(83, 10).Value = 0
$ws.Cells.Item(83, 12).Value = 0
$ws.Cells.Item(83, 14).ClearContents()

$ws.Cells.Item(84, 8).Value = 100017300
$ws.Cells.Item(84, 10).Value = 100017300
$ws.Cells.Item(84, 12).Value = 900155700
$ws.Cells.Item(84, 14).Value = -900166932

$ws.Cells.Item(86, 8).Value = 5000
$ws.Cells.Item(86, 9).Value = 0
$ws.Cells.Item(86, 11).Value = 0
$ws.Cells.Item(86, 13).ClearContents()

$ws.Cells.Item(88, 8).Value = 10999.6
$ws.Cells.Item(88, 10).Value = 10999.6
$ws.Cells.Item(88, 12).Value = 32998.8
$ws.Cells.Item(88, 14).Value = -33854.8

$ws.Cells.Item(89, 8).Value = 5000
$ws.Cells.Item(89, 9).Value = 0
$ws.Cells.Item(89, 11).Value = 0
$ws.Cells.Item(89, 13).ClearContents()

$ws.Cells.Item(91, 8).Value = 10999.6
$ws.Cells.Item(91, 10).Value = 10999.6
$ws.Cells.Item(91, 12).Value = 32998.8
$ws.Cells.Item(91, 14).Value = -35962.8

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 5660.6665
$ws.Cells.Item(70, 9).Value = 2000
$ws.Cells.Item(70, 10).Value = 6392.8
$ws.Cells.Item(70, 11).Value = 2000
$ws.Cells.Item(70, 12).Value = 6392.8
$ws.Cells.Item(70, 13).Value = -1730
$ws.Cells.Item(70, 14).Value = -6932.8

$ws.Cells.Item(73, 8).Value = 5660.6665
$ws.Cells.Item(73, 9).Value = 2000
$ws.Cells.Item(73, 10).Value = 6392.8
$ws.Cells.Item(73, 11).Value = 2000
$ws.Cells.Item(73, 12).Value = 6392.8
$ws.Cells.Item(73, 13).Value = -1064
$ws.Cells.Item(73, 14).Value = -8264.799999999999

$ws.Cells.Item(132, 8).Value = 1972.5641
$ws.Cells.Item(132, 9).Value = 1748.6285
$ws.Cells.Item(132, 10).Value = 3932
$ws.Cells.Item(132, 11).Value = 5245.8855
$ws.Cells.Item(132, 12).Value = 11796
$ws.Cells.Item(132, 13).Value = -2715.8855
$ws.Cells.Item(132, 14).Value = -16856

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(16, 8).Value = 1184.2858
$ws.Cells.Item(16, 9).Value = 1400.7
$ws.Cells.Item(16, 10).Value = 643.25
$ws.Cells.Item(16, 11).Value = 1400.7
$ws.Cells.Item(16, 12).Value = 643.25
$ws.Cells.Item(16, 13).Value = -1230.7
$ws.Cells.Item(16, 14).Value = -983.25

$ws.Cells.Item(36, 8).Value = 41999.5
$ws.Cells.Item(36, 10).Value = 41999.5
$ws.Cells.Item(36, 12).Value = 41999.5
$ws.Cells.Item(36, 14).Value = -43123.5

$ws.Cells.Item(68, 8).Value = 2363.8333
$ws.Cells.Item(68, 9).Value = 1932.9166
$ws.Cells.Item(68, 11).Value = 1932.9166
$ws.Cells.Item(68, 13).Value = -1183.9166

$ws.Cells.Item(71, 8).Value = 2363.8333
$ws.Cells.Item(71, 9).Value = 1932.9166
$ws.Cells.Item(71, 11).Value = 9664.583000000001
$ws.Cells.Item(71, 13).Value = -5920.583000000001

$ws.Cells.Item(124, 8).Value = 69999
$ws.Cells.Item(124, 10).Value = 69999
$ws.Cells.Item(124, 12).Value = 69999
$ws.Cells.Item(124, 14).Value = -79819

$ws.Cells.Item(133, 8).Value = 90225
$ws.Cells.Item(133, 10).Value = 90225
$ws.Cells.Item(133, 12).Value = 90225
$ws.Cells.Item(133, 14).Value = -95285

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(136, 8).Value = 2391.6191
$ws.Cells.Item(136, 9).Value = 2475.25
$ws.Cells.Item(136, 11).Value = 7425.75
$ws.Cells.Item(136, 13).Value = -4875.75
